$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 37049300
$ws.Range("E8").Value = 37362100
$ws.Range("F8").Value = 42843200
$ws.Range("G8").Value = 42969000
$ws.Range("H8").Value = 43052500
$ws.Range("I8").Value = 39610800
$ws.Range("J8").Value = 40386900
$ws.Range("D9").Value = 26817900
$ws.Range("E9").Value = 26831900
$ws.Range("F9").Value = 31530000
$ws.Range("G9").Value = 31384300
$ws.Range("H9").Value = 31584400
$ws.Range("I9").Value = 28728800
$ws.Range("J9").Value = 29218600
$ws.Range("D10").Value = 10231400
$ws.Range("E10").Value = 10530100
$ws.Range("F10").Value = 11313200
$ws.Range("G10").Value = 11584800
$ws.Range("H10").Value = 11468100
$ws.Range("I10").Value = 10882000
$ws.Range("J10").Value = 11168300
$ws.Range("H12").Value = 2001400
$ws.Range("I12").Value = 2088700
$ws.Range("J12").Value = 2154800
$ws.Range("E14").Value = 380500
$ws.Range("H14").Value = 568500
$ws.Range("I14").Value = 1378900
$ws.Range("J14").Value = 248400
$ws.Range("D17").Value = 35399600
$ws.Range("E17").Value = 36300300
$ws.Range("F17").Value = 41752900
$ws.Range("G17").Value = 41354200
$ws.Range("H17").Value = 41721100
$ws.Range("I17").Value = 40191700
$ws.Range("J17").Value = 39683300
$ws.Range("D18").Value = 1649700
$ws.Range("E18").Value = 1061800
$ws.Range("F18").Value = 1090300
$ws.Range("G18").Value = 1614800
$ws.Range("H18").Value = 1331400
$ws.Range("I18").Value = -580900
$ws.Range("J18").Value = 703500
$ws.Range("D20").Value = 578600
$ws.Range("E20").Value = 99100
$ws.Range("F20").Value = 148800
$ws.Range("G20").Value = 236900
$ws.Range("H20").Value = 184200
$ws.Range("I20").Value = 175600
$ws.Range("J20").Value = -16500
$ws.Range("D21").Value = 3764600
$ws.Range("E21").Value = 2864900
$ws.Range("F21").Value = 3077300
$ws.Range("G21").Value = 3551700
$ws.Range("H21").Value = 3290800
$ws.Range("I21").Value = 1337900
$ws.Range("J21").Value = 2583400
$ws.Range("D22").Value = 36200
$ws.Range("E22").Value = 38400
$ws.Range("F22").Value = 47400
$ws.Range("G22").Value = 54000
$ws.Range("H22").Value = 59200
$ws.Range("I22").Value = 65900
$ws.Range("J22").Value = 83900
$ws.Range("D23").Value = 2192100
$ws.Range("E23").Value = 1122400
$ws.Range("F23").Value = 1191700
$ws.Range("G23").Value = 1797700
$ws.Range("H23").Value = 1456400
$ws.Range("I23").Value = -471200
$ws.Range("J23").Value = 603100
$ws.Range("D24").Value = 673200
$ws.Range("E24").Value = 282000
$ws.Range("F24").Value = 374300
$ws.Range("G24").Value = 486800
$ws.Range("H24").Value = 353400
$ws.Range("I24").Value = 219300
$ws.Range("J24").Value = 271200
$ws.Range("D26").Value = 1518900
$ws.Range("E26").Value = 840400
$ws.Range("F26").Value = 817400
$ws.Range("G26").Value = 1310900
$ws.Range("H26").Value = 1103000
$ws.Range("I26").Value = -690500
$ws.Range("J26").Value = 331900
$ws.Range("D27").Value = 1447400
$ws.Range("E27").Value = 778700
$ws.Range("F27").Value = 784300
$ws.Range("G27").Value = 1265800
$ws.Range("H27").Value = 957800
$ws.Range("I27").Value = -722500
$ws.Range("J27").Value = 386100
$ws.Range("D29").Value = 83400
$ws.Range("E29").Value = 21300
$ws.Range("D32").Value = -578600
$ws.Range("E32").Value = -99100
$ws.Range("F32").Value = -148800
$ws.Range("G32").Value = -236900
$ws.Range("H32").Value = -184200
$ws.Range("I32").Value = -175600
$ws.Range("J32").Value = 16500
$ws.Range("D33").Value = 1530800
$ws.Range("E33").Value = 799900
$ws.Range("F33").Value = 784300
$ws.Range("G33").Value = 1265800
$ws.Range("H33").Value = 957800
$ws.Range("I33").Value = -722500
$ws.Range("J33").Value = 386100
$ws.Range("D35").Value = 1530800
$ws.Range("E35").Value = 799900
$ws.Range("F35").Value = 784300
$ws.Range("G35").Value = 1265800
$ws.Range("H35").Value = 957800
$ws.Range("I35").Value = -722500
$ws.Range("J35").Value = 386100
$ws.Range("D41").Value = 4091100
$ws.Range("E41").Value = 3441500
$ws.Range("F41").Value = 3442500
$ws.Range("G41").Value = 3272700
$ws.Range("H41").Value = 4843200
$ws.Range("I41").Value = 1830600
$ws.Range("J41").Value = 1930000
$ws.Range("H42").Value = 642400
$ws.Range("I42").Value = 926300
$ws.Range("J42").Value = 546300
$ws.Range("D43").Value = 9504100
$ws.Range("E43").Value = 9637400
$ws.Range("F43").Value = 9750500
$ws.Range("G43").Value = 10213000
$ws.Range("H43").Value = 18675900
$ws.Range("I43").Value = 7990500
$ws.Range("J43").Value = 8032200
$ws.Range("D44").Value = 2184100
$ws.Range("E44").Value = 2650200
$ws.Range("F44").Value = 2701600
$ws.Range("G44").Value = 2837500
$ws.Range("H44").Value = 4172200
$ws.Range("I44").Value = 2920800
$ws.Range("J44").Value = 3020400
$ws.Range("D45").Value = 1101100
$ws.Range("E45").Value = 926200
$ws.Range("F45").Value = 773900
$ws.Range("G45").Value = 741400
$ws.Range("H45").Value = 2879100
$ws.Range("I45").Value = 1901400
$ws.Range("J45").Value = 1855300
$ws.Range("D46").Value = 16880400
$ws.Range("E46").Value = 16655300
$ws.Range("F46").Value = 16668500
$ws.Range("G46").Value = 17064500
$ws.Range("H46").Value = 16137800
$ws.Range("I46").Value = 15569500
$ws.Range("J46").Value = 15384100
$ws.Range("D47").Value = 2818600
$ws.Range("E47").Value = 2638100
$ws.Range("F47").Value = 2359800
$ws.Range("G47").Value = 2454200
$ws.Range("H47").Value = 3522000
$ws.Range("I47").Value = 3086900
$ws.Range("J47").Value = 1347800
$ws.Range("D48").Value = 4751300
$ws.Range("E48").Value = 5393700
$ws.Range("F48").Value = 5572900
$ws.Range("G48").Value = 5744800
$ws.Range("H48").Value = 11228600
$ws.Range("I48").Value = 5590900
$ws.Range("J48").Value = 5794100
$ws.Range("D49").Value = 1565500
$ws.Range("E49").Value = 1764700
$ws.Range("F49").Value = 1849100
$ws.Range("G49").Value = 1854800
$ws.Range("H49").Value = 3459000
$ws.Range("I49").Value = 1693400
$ws.Range("J49").Value = 2081800
$ws.Range("D52").Value = 2202800
$ws.Range("E52").Value = 2399300
$ws.Range("F52").Value = 2715500
$ws.Range("G52").Value = 2452600
$ws.Range("H52").Value = 3750000
$ws.Range("I52").Value = 4104900
$ws.Range("J52").Value = 2019500
$ws.Range("D54").Value = 28218600
$ws.Range("E54").Value = 28851100
$ws.Range("F54").Value = 29165800
$ws.Range("G54").Value = 29570900
$ws.Range("H54").Value = 28077700
$ws.Range("I54").Value = 26399700
$ws.Range("J54").Value = 26627400
$ws.Range("D57").Value = 4699400
$ws.Range("E57").Value = 5352500
$ws.Range("F57").Value = 5516200
$ws.Range("G57").Value = 5967600
$ws.Range("H57").Value = 11581200
$ws.Range("I57").Value = 5123500
$ws.Range("J57").Value = 5584500
$ws.Range("D58").Value = 1223900
$ws.Range("E58").Value = 1182300
$ws.Range("F58").Value = 1307800
$ws.Range("G58").Value = 1558500
$ws.Range("H58").Value = 2578500
$ws.Range("I58").Value = 5185600
$ws.Range("J58").Value = 1308200
$ws.Range("D59").Value = 6033300
$ws.Range("E59").Value = 6410400
$ws.Range("F59").Value = 6257300
$ws.Range("G59").Value = 6244800
$ws.Range("H59").Value = 9398300
$ws.Range("I59").Value = 6306700
$ws.Range("J59").Value = 5920900
$ws.Range("D60").Value = 11956700
$ws.Range("E60").Value = 12945200
$ws.Range("F60").Value = 13081300
$ws.Range("G60").Value = 13770900
$ws.Range("H60").Value = 13410500
$ws.Range("I60").Value = 14179300
$ws.Range("J60").Value = 12813600
$ws.Range("D61").Value = 2409200
$ws.Range("E61").Value = 3202900
$ws.Range("F61").Value = 3506500
$ws.Range("G61").Value = 3671000
$ws.Range("H61").Value = 3775400
$ws.Range("I61").Value = 2459000
$ws.Range("J61").Value = 2530800
$ws.Range("D62").Value = 2960400
$ws.Range("E62").Value = 3489400
$ws.Range("F62").Value = 4204800
$ws.Range("G62").Value = 3682000
$ws.Range("H62").Value = 9071500
$ws.Range("I62").Value = 4572900
$ws.Range("J62").Value = 2544900
$ws.Range("D66").Value = 18384900
$ws.Range("E66").Value = 20884300
$ws.Range("F66").Value = 22089400
$ws.Range("G66").Value = 22428500
$ws.Range("H66").Value = 22956400
$ws.Range("I66").Value = 20757700
$ws.Range("J66").Value = 19023700
$ws.Range("D72").Value = 4337200
$ws.Range("E72").Value = 2403700
$ws.Range("F72").Value = 1409600
$ws.Range("G72").Value = 1181900
$ws.Range("H72").Value = 2344400
$ws.Range("I72").Value = 2395000
$ws.Range("J72").Value = 3302300
$ws.Range("D76").Value = 9833700
$ws.Range("E76").Value = 7966900
$ws.Range("F76").Value = 7076300
$ws.Range("G76").Value = 7142400
$ws.Range("H76").Value = 5121300
$ws.Range("I76").Value = 5642100
$ws.Range("J76").Value = 7603700
$ws.Range("D81").Value = 1530800
$ws.Range("E81").Value = 799900
$ws.Range("F81").Value = 784300
$ws.Range("G81").Value = 1265800
$ws.Range("H81").Value = 957800
$ws.Range("I81").Value = -722500
$ws.Range("J81").Value = 386100
$ws.Range("D83").Value = 1534500
$ws.Range("E83").Value = 1702200
$ws.Range("F83").Value = 1836200
$ws.Range("G83").Value = 1698100
$ws.Range("H83").Value = 1773200
$ws.Range("I83").Value = 1741300
$ws.Range("J83").Value = 1894300
$ws.Range("D89").Value = 1811800
$ws.Range("E89").Value = 2263000
$ws.Range("F89").Value = 2288000
$ws.Range("G89").Value = 2532500
$ws.Range("H89").Value = 1595600
$ws.Range("I89").Value = 646700
$ws.Range("J89").Value = 2169700
$ws.Range("D91").Value = -1220400
$ws.Range("E91").Value = -1794400
$ws.Range("F91").Value = -1715500
$ws.Range("G91").Value = -1799900
$ws.Range("H91").Value = -2655900
$ws.Range("I91").Value = -1008200
$ws.Range("J91").Value = -1245600
$ws.Range("D94").Value = -204100
$ws.Range("E94").Value = -1315100
$ws.Range("F94").Value = -1485400
$ws.Range("G94").Value = -1812700
$ws.Range("H94").Value = -1165600
$ws.Range("I94").Value = -1459800
$ws.Range("J94").Value = -1725100
$ws.Range("D96").Value = -185400
$ws.Range("E96").Value = -149600
$ws.Range("F96").Value = -149600
$ws.Range("G96").Value = -149600
$ws.Range("H96").Value = -24100
$ws.Range("I96").Value = -208900
$ws.Range("J96").Value = -204900
$ws.Range("D100").Value = -1017000
$ws.Range("E100").Value = -894000
$ws.Range("F100").Value = -612400
$ws.Range("G100").Value = -156600
$ws.Range("H100").Value = -417800
$ws.Range("I100").Value = 907500
$ws.Range("J100").Value = -1256300
$ws.Range("D101").Value = 30400
$ws.Range("E101").Value = -25300
$ws.Range("F101").Value = -20400
$ws.Range("G101").Value = -13000
$ws.Range("H101").Value = 138000
$ws.Range("I101").Value = 67000
$ws.Range("J101").Value = -56100
$ws.Range("D102").Value = 621100
$ws.Range("E102").Value = 28600
$ws.Range("F102").Value = 169800
$ws.Range("G102").Value = 550200
$ws.Range("H102").Value = 150200
$ws.Range("I102").Value = 161400
$ws.Range("J102").Value = -867800
